$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PacketFormat")

# Insert a new row at 19, shifting the "Environmental" (0x55) and "Multitemp" (0x56)
# blocks down by one row.
$ws.Rows("19:19").Insert()

# New row 19: CounterV2 / 0x54 extension type header.
$ws.Range("L19").Value = "CounterV2"
$ws.Range("M19").Value = "0x54"

# Row 21 (former Multitemp row, now shifted down): extend to 4 temperature
# values (hi/lo pairs) instead of the old 7:0 nibble layout, and add two more
# columns (T21:U21) for symmetry with the header row below.
$ws.Range("N21").Value = "Temp1[hi]"
$ws.Range("O21").Value = "Temp1[lo]"
$ws.Range("Q21").Value = "Temp2[lo]"
$ws.Range("P21").Value = "Temp2[hi]"
$ws.Range("S21").Value = "Temp3[lo]"
$ws.Range("R21").Value = "Temp3[hi]"
$ws.Range("S21").Copy()
$ws.Range("T21").PasteSpecial(-4122)
$ws.Range("U21").PasteSpecial(-4122)

# Row 16 (index row above L17:S20 table) gains two more running indices.
$ws.Range("S16").Copy()
$ws.Range("T16").PasteSpecial(-4122)
$ws.Range("U16").PasteSpecial(-4122)
$ws.Range("T16").Value = 16
$ws.Range("U16").Value = 17

$ws.Range("T21:U21").Select()
